$wb = $excel.ActiveWorkbook

# Rename the "wt" and "dcin5" sheets to reflect log2 expression data
$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Update the selection on the wt_log2_expression sheet
$wsWt.Range("C38").Select()

# Restore the originally active sheet so the workbook's active tab is unchanged
$wb.Worksheets.Item("network_optimized_weights").Activate()
